# Added in fancy stats boxes and deployed server
#
# 1. Students sheet: Eliezer Yudkowsky redeemed a code (Codes count 0 -> 1)
# 2. Books sheet: add "Rationality: from AI to Zombies" by Eliezer Yudkowsky
#    as a new row, inserted alphabetically before "Smart People Should Build Things"
# 3. Codes sheet: add the new book's redeemed code (Yudkowsky, Eliezer) plus a
#    fresh unredeemed code for "Life of Pi"

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Students -------------------------------------------------
$students = $wb.Worksheets.Item("Students")
$students.Range("E8").Value = 1

# --- Sheet 2: Books ------------------------------------------------------
$books = $wb.Worksheets.Item("Books")

# Insert a new row 4 (pushes "Smart People Should Build Things" and the
# rows below it down by one) and fill in the new book's data.
$books.Rows.Item(4).Insert()
$books.Range("A4").Value = "Rationality: from AI to Zombies"
$books.Range("B4").Value = "Eliezer Yudkowsky"
$books.Range("C4").Value = 1
$books.Range("D4").Value = 1

# --- Sheet 3: Codes --------------------------------------------------------
$codes = $wb.Worksheets.Item("Codes")

# Append two new rows after the existing last row (row 5), inheriting the
# formatting of the row above via Insert().
$codes.Rows.Item(6).Insert()
$codes.Range("A6").Value = "Rationality: from AI to Zombies"
$codes.Range("B6").Value = "W7CPzamcGj"
$codes.Range("C6").Value = "Yudkowsky, Eliezer"

$codes.Rows.Item(7).Insert()
$codes.Range("A7").Value = "Life of Pi"
$codes.Range("B7").Value = "we1zhU2xyb"
$codes.Range("C7").Value = "None"
